$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they remain strings like the original.
$textCells = @("D5", "D6", "D10", "D14", "D15", "D17", "D19", "D20", "D21", "D24", "D26", "D27", "D29", "D30", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the latest cryptos data refresh
$ws.Range("D2").Value = "55.314.57"
$ws.Range("E2").Value = "  -4.24%  "
$ws.Range("D3").Value = "2.904.75"
$ws.Range("E3").Value = "  -4.26%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "500.10"
$ws.Range("E5").Value = "  -2.33%  "
$ws.Range("D6").Value = "132.12"
$ws.Range("E6").Value = "  -5.70%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -5.36%  "
$ws.Range("E9").Value = "  -4.79%  "
$ws.Range("D10").Value = "0.103"
$ws.Range("E10").Value = "  -6.62%  "
$ws.Range("E11").Value = "  -5.64%  "
$ws.Range("D12").Value = "3.406.57"
$ws.Range("E12").Value = "  -4.05%  "
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").Value = "25.60"
$ws.Range("E14").Value = "  -4.44%  "
$ws.Range("D15").Value = "0.0000157"
$ws.Range("E15").Value = "  -5.65%  "
$ws.Range("D16").Value = "55.344.88"
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("D17").Value = "5.96"
$ws.Range("E17").Value = "  -4.65%  "
$ws.Range("D18").Value = "2.903.97"
$ws.Range("E18").Value = "  -4.30%  "
$ws.Range("D19").Value = "12.46"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").Value = "7.63"
$ws.Range("E20").Value = "  -4.90%  "
$ws.Range("D21").Value = "311.50"
$ws.Range("E21").Value = "  -6.74%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  -3.91%  "
$ws.Range("D24").Value = "62.44"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").Value = "3.029.00"
$ws.Range("E25").Value = "  -3.98%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "0.158"
$ws.Range("E27").Value = "  -6.90%  "
$ws.Range("D28").Value = "0.0₃0829"
$ws.Range("E28").Value = "  -11.07%  "
$ws.Range("D29").Value = "6.32"
$ws.Range("E29").Value = "  -7.82%  "
$ws.Range("D30").Value = "6.88"
$ws.Range("E30").Value = "  -8.18%  "
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").Value = "19.69"
$ws.Range("E32").Value = "  -5.61%  "
$ws.Range("E33").Value = "  -8.56%  "
$ws.Range("D34").Value = "149.46"
$ws.Range("E34").Value = "  -4.14%  "
$ws.Range("D35").Value = "4.34"
$ws.Range("E35").Value = "  -8.17%  "
$ws.Range("D36").Value = "5.54"
$ws.Range("E36").Value = "  -5.76%  "
$ws.Range("D37").Value = "24.09"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").Value = "1.17"
$ws.Range("E38").Value = "  -8.63%  "
$ws.Range("D39").Value = "0.0641"
$ws.Range("E39").Value = "  -6.82%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "36.14"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("D42").Value = "3.67"
$ws.Range("E42").Value = "  -5.56%  "
$ws.Range("D43").Value = "0.633"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").Value = "2.086.49"
$ws.Range("E44").Value = "  -9.82%  "
$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").Value = "5.87"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.31"
$ws.Range("E46").Value = "  -8.10%  "
$ws.Range("D47").Value = "0.903"
$ws.Range("E47").Value = "  -9.37%  "
$ws.Range("D48").Value = "0.0230"
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("D49").Value = "18.50"
$ws.Range("E49").Value = "  -6.80%  "
$ws.Range("D50").Value = "0.0832"
$ws.Range("E50").Value = "  -7.11%  "
$ws.Range("E51").Value = "  -10.21%  "
